# Sprint Burndown Chart Template - update
#
# Commit: "Updated SprintBurnDown and Finished LOG_IN Scenario from Salonika"
#
# The substantive edit is a single data correction on the burndown tracking
# sheet: the "Added" effort for the 3rd sprint day (cell E18) is corrected
# from 1 to 6. Everything else on the sheet (the running totals in row 26,
# the ideal-trend projection in row 27, and the burndown chart series that
# plot those two rows) is formula-driven off of this input, so updating the
# single cell lets Excel's recalculation engine ripple the new numbers
# through the rest of the workbook automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Added" effort entry for the third day of the sprint.
$ws.Range("E18").Value = 6

# Leave the cursor where the author left it when they saved the file.
$ws.Range("G19").Select()
